# Apply the "Upload new version with timestamp" edit:
#  - AUGRAM 1GM 14 TABS: ratio column (H) changes from 1:0 -> 0:1
#  - Insert a new item "GLIPTUS PLUS 50/1000MG 30 TABLETS" after CETAL COLD & FLU (new row 10)
#  - Append a new item "VITACID C 1GM 12 EFF. TAB." after SIDERAL ACTIVE (new row 15)
#  - Refresh the total (sum of sell-price column) and the generated timestamp footer

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) AUGRAM row (row 8): ratio value changes
# ---------------------------------------------------------------------------
$ws.Range("H8").Value = "0:1"

# ---------------------------------------------------------------------------
# 2) Insert a brand-new data row (10) for GLIPTUS PLUS, pushing the existing
#    rows 10-15 (METFORMIN..footer) down by one.
# ---------------------------------------------------------------------------
$ws.Rows(10).Insert()

# Bring over formatting (styles) for the new row from the row above it, then
# fix up the row height / merges to match the rest of the table.
$ws.Range("A9:Q9").Copy($ws.Range("A10:Q10"))
$ws.Rows(10).RowHeight = 24.75
$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "GLIPTUS PLUS 50/1000MG 30 TABLETS"
$ws.Range("H10").Value = "0:0"

# L/P columns carry a numeric-looking NumberFormat, so force text storage
# (matching the source file, which stores these as shared-string text)
# by flipping to a text format, assigning, then restoring the format.
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "1"
$ws.Range("L10").NumberFormat = '#,##0.##;"["#,##0.##"]";0'

$ws.Range("N10").Value = "192.00"

$ws.Range("P10").NumberFormat = "@"
$ws.Range("P10").Value = "192.0000"
$ws.Range("P10").NumberFormat = "0.00"

$ws.Range("Q10").Value = "1:0"

# Renumber the rows that followed (old #4-#7 -> new #5-#8)
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6
$ws.Range("A13").Value = 7
$ws.Range("A14").Value = 8

# ---------------------------------------------------------------------------
# 3) Append two new data rows (15) for VITACID, after SIDERAL ACTIVE (row 14),
#    pushing the totals row (15) and footer row (16) down by one.
# ---------------------------------------------------------------------------
$ws.Rows(15).Insert()

$ws.Range("A14:Q14").Copy($ws.Range("A15:Q15"))
$ws.Rows(15).RowHeight = 24.75
$ws.Range("A15:B15").Merge()
$ws.Range("C15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()
$ws.Range("N15:O15").Merge()

$ws.Range("A15").Value = 9
$ws.Range("C15").Value = "VITACID C 1GM 12 EFF. TAB."
$ws.Range("H15").Value = "0:0"
$ws.Range("L15").Value = "1"
$ws.Range("N15").Value = "54.00"
$ws.Range("P15").Value = "54.0000"
$ws.Range("Q15").Value = "1:0"

# ---------------------------------------------------------------------------
# 4) Refresh the total (now row 16) and footer timestamp (now row 17)
# ---------------------------------------------------------------------------
$ws.Range("P16").Value = 655.8
$ws.Range("A17").Value = "Saturday, 16 August, 2025 9:16 PM"
